# Applies the "cryptos list" data refresh described in the commit:
# "Updated cryptos list on Mon Mar 27 15:44:04 UTC 2023 with GitHub Actions"
#
# The sheet stores Price (column D) and Volume(1h) (column E) as text cells
# (t="inlineStr" in the original file) even though many Price values look like
# plain numbers. When Excel is given a numeric-looking string it auto-converts
# it to a real number (dropping e.g. a trailing "0"), so for any new Price value
# that parses as a number we force the cell to text (NumberFormat "@") before
# assigning it, then restore the cell style to "Normal" so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.009.45"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "1.717.41"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3432"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07289"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.050"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.863"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "1.728.25"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06286"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.622"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "27.088.36"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("E24").Value = "  -4.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "1.926.12"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.143"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.021"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09074"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.594"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02195"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05819"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.771"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("E39").Value = "  -5.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.402"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5964"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.478"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.65%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.03%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.616"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5615"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.860"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.088"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
